$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.195.56"
$ws.Range("E2").Value = "  +4.18%  "
$ws.Range("D3").Value = "2.627.89"
$ws.Range("E3").Value = "  +4.65%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.65%  "
$ws.Range("D9").Value = "2.625.94"
$ws.Range("E9").Value = "  +4.58%  "
$ws.Range("E10").Value = "  +15.34%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("D14").Value = "3.081.87"
$ws.Range("E14").Value = "  +3.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("E16").Value = "  +7.80%  "
$ws.Range("D17").Value = "71.182.51"
$ws.Range("E17").Value = "  +4.50%  "
$ws.Range("D18").Value = "2.622.91"
$ws.Range("E18").Value = "  +4.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "382.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.96%  "
$ws.Range("E20").Value = "  +5.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("E24").Value = "  +5.85%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.92%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.46%  "
$ws.Range("D28").Value = "2.759.75"
$ws.Range("E28").Value = "  +4.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.92%  "
$ws.Range("D30").Value = "0.0₃0961"
$ws.Range("E30").Value = "  +7.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "539.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("E33").Value = "  +4.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.41%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.66%  "
$ws.Range("E39").Value = "  +7.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.38"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.32%  "
$ws.Range("E42").Value = "  +9.62%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.69%  "
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "154.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0267"
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.533"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.51%  "
